$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to mimic the "%.16g" float formatting used by the original authoring
# tool (openpyxl), so recomputed numbers land on the exact same IEEE754
# double as the target file once the decimal text is re-parsed.
function Round16([double]$x) {
    return [double]($x.ToString("G16"))
}

$lastRow = 81

# Step 1: Insert a new column before column B ("pixel"), shifting
# pixel -> C, volume -> D, jarak -> E.
$ws.Columns.Item(2).Insert()

# Step 2: Fix up formatting. The Insert() operation carried column A's
# header/border style into the new column B, which is wrong: the header
# cell B1 should look like the other header cells (style copied from C1),
# while the data cells B2:B81 should look like plain data cells (style
# copied from C2:C81, i.e. no special border/bold).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C2:C" + $lastRow).Copy()
$ws.Range("B2:B" + $lastRow).PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Step 3: Write the new header label.
$ws.Range("B1").Value = "Unnamed: 0"

# Step 4: Populate the new column B with the same running index values
# already present in column A.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value()
}

# Step 5: Rescale the "volume" column (now column D) from its old unit to
# the new one (multiply by 1,000,000).
for ($r = 2; $r -le $lastRow; $r++) {
    $old = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 4).Value = Round16($old * 1000000)
}
